$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" footer date field from 12/4/2023 to
#    5/19/2025 everywhere it is defined: the slide master, every one of its
#    custom (slide) layouts, and the notes master.
# ---------------------------------------------------------------------------
$newDate = "5/19/2025"

function Set-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster

# slide master's own date placeholder
Set-DatePlaceholderText $master.Shapes

# every custom layout hanging off the slide master
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Set-DatePlaceholderText $layouts.Item($L).Shapes
}

# notes master's date placeholder
Set-DatePlaceholderText $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 1 - "Rectangle 70": drop the second "YOUSIF ALOUFI" paragraph,
#    keep "JANI SHARIFF SHAIK" only (box auto-shrinks to the single line).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 70") {
        $shp.TextFrame.TextRange.Text = "JANI SHARIFF SHAIK"
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 22 - "Rectangle 4": drop the second "YOUSIF SAEED ALOUFI"
#    paragraph, keep "JANI SHARIFF SHAIK" only, and reposition the box.
# ---------------------------------------------------------------------------
$slide22 = $p.Slides.Item(22)
for ($i = 1; $i -le $slide22.Shapes.Count; $i++) {
    $shp = $slide22.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 4") {
        $shp.TextFrame.TextRange.Text = "JANI SHARIFF SHAIK"
        $shp.Left = 8073632 / 914400 * 72
        $shp.Top = 6073171 / 914400 * 72
    }
}
